$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($rng, $val)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "34.545.25"
Set-TextValue $ws.Range("E2") "  +1.73%  "
Set-TextValue $ws.Range("D3") "1.843.56"
Set-TextValue $ws.Range("E3") "  +4.10%  "
Set-TextValue $ws.Range("D4") "0.999"
Set-TextValue $ws.Range("E4") "  -0.19%  "
Set-TextValue $ws.Range("D5") "226.29"
Set-TextValue $ws.Range("E5") "  +0.61%  "
Set-TextValue $ws.Range("D6") "0.555"
Set-TextValue $ws.Range("E6") "  +1.62%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.21%  "
Set-TextValue $ws.Range("D8") "32.45"
Set-TextValue $ws.Range("E8") "  +4.44%  "
Set-TextValue $ws.Range("D9") "0.294"
Set-TextValue $ws.Range("E9") "  +5.79%  "
Set-TextValue $ws.Range("D10") "0.0718"
Set-TextValue $ws.Range("E10") "  +10.01%  "
Set-TextValue $ws.Range("D11") "0.0933"
Set-TextValue $ws.Range("E11") "  +0.62%  "
Set-TextValue $ws.Range("D12") "2.113.38"
Set-TextValue $ws.Range("E12") "  +4.46%  "
Set-TextValue $ws.Range("D13") "1.847.93"
Set-TextValue $ws.Range("E13") "  +4.29%  "
Set-TextValue $ws.Range("D14") "11.11"
Set-TextValue $ws.Range("E14") "  +1.95%  "
Set-TextValue $ws.Range("D15") "0.652"
Set-TextValue $ws.Range("E15") "  +5.37%  "
Set-TextValue $ws.Range("D16") "34.568.33"
Set-TextValue $ws.Range("E16") "  +1.86%  "
Set-TextValue $ws.Range("D17") "4.36"
Set-TextValue $ws.Range("E17") "  +4.46%  "
Set-TextValue $ws.Range("D18") "69.90"
Set-TextValue $ws.Range("E18") "  +2.31%  "
Set-TextValue $ws.Range("D19") "253.18"
Set-TextValue $ws.Range("E19") "  +1.13%  "
Set-TextValue $ws.Range("D20") "0.0₃0807"
Set-TextValue $ws.Range("E20") "  +10.15%  "
Set-TextValue $ws.Range("D21") "11.32"
Set-TextValue $ws.Range("E21") "  +10.39%  "
Set-TextValue $ws.Range("D22") "0.997"
Set-TextValue $ws.Range("E22") "  -0.52%  "
Set-TextValue $ws.Range("D23") "4.33"
Set-TextValue $ws.Range("E23") "  +3.82%  "
Set-TextValue $ws.Range("D24") "2.15"
Set-TextValue $ws.Range("E24") "  +1.25%  "
Set-TextValue $ws.Range("D25") "162.01"
Set-TextValue $ws.Range("E25") "  +4.44%  "
Set-TextValue $ws.Range("D26") "16.89"
Set-TextValue $ws.Range("E26") "  +3.63%  "
Set-TextValue $ws.Range("D27") "7.27"
Set-TextValue $ws.Range("E27") "  +4.60%  "
Set-TextValue $ws.Range("D28") "0.115"
Set-TextValue $ws.Range("E28") "  +2.17%  "
Set-TextValue $ws.Range("D29") "0.998"
Set-TextValue $ws.Range("E29") "  -0.26%  "
Set-TextValue $ws.Range("D30") "0.0539"
Set-TextValue $ws.Range("E30") "  +5.89%  "
Set-TextValue $ws.Range("D31") "3.82"
Set-TextValue $ws.Range("E31") "  +2.23%  "
Set-TextValue $ws.Range("D32") "1.21"
Set-TextValue $ws.Range("E32") "  +1.96%  "
Set-TextValue $ws.Range("D33") "508.55"
Set-TextValue $ws.Range("E33") "  +873.32%  "
Set-TextValue $ws.Range("D34") "3.65"
Set-TextValue $ws.Range("E34") "  +3.10%  "
Set-TextValue $ws.Range("D35") "1.95"
Set-TextValue $ws.Range("E35") "  +7.23%  "
Set-TextValue $ws.Range("D36") "1.456.52"
Set-TextValue $ws.Range("E36") "  +1.00%  "
Set-TextValue $ws.Range("D37") "0.657"
Set-TextValue $ws.Range("E37") "  +6.18%  "
Set-TextValue $ws.Range("E38") "  +2.85%  "
Set-TextValue $ws.Range("E39") "  +5.55%  "
Set-TextValue $ws.Range("D40") "0.979"
Set-TextValue $ws.Range("E40") "  +11.39%  "
Set-TextValue $ws.Range("B41") "Aave"
Set-TextValue $ws.Range("C41") "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-TextValue $ws.Range("D41") "83.10"
Set-TextValue $ws.Range("B42") "MXToken"
Set-TextValue $ws.Range("C42") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D42") "2.81"
Set-TextValue $ws.Range("E42") "  -1.25%  "
Set-TextValue $ws.Range("D43") "2.38"
Set-TextValue $ws.Range("E43") "  +1.08%  "
Set-TextValue $ws.Range("D44") "2.15"
Set-TextValue $ws.Range("E44") "  +5.97%  "
Set-TextValue $ws.Range("D45") "6.13"
Set-TextValue $ws.Range("E45") "  +7.52%  "
Set-TextValue $ws.Range("D46") "2.011.06"
Set-TextValue $ws.Range("E46") "  +4.68%  "
Set-TextValue $ws.Range("B47") "WEMIXToken"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
Set-TextValue $ws.Range("D47") "1.06"
Set-TextValue $ws.Range("E47") "  +0.85%  "
Set-TextValue $ws.Range("B48") "Kaspa"
Set-TextValue $ws.Range("C48") "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
Set-TextValue $ws.Range("D48") "0.0497"
Set-TextValue $ws.Range("E48") "  -2.04%  "
Set-TextValue $ws.Range("B49") "InjectiveProtocol"
Set-TextValue $ws.Range("C49") "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue $ws.Range("D49") "12.34"
Set-TextValue $ws.Range("E49") "  +5.05%  "
Set-TextValue $ws.Range("D50") "106.71"
Set-TextValue $ws.Range("E50") "  +10.36%  "
Set-TextValue $ws.Range("E51") "  +0.09%  "
